$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7:G7").NumberFormat = "@"

$ws.Range("A7").Value = "address f6"
$ws.Range("B7").Value = "f6 city"
$ws.Range("C7").Value = "first name f6"
$ws.Range("D7").Value = "ls name f6"
$ws.Range("E7").Value = " "
$ws.Range("F7").Value = "4444444"
$ws.Range("G7").Value = "4444666"
